$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.19"
$ws.Range("E2").Value = "'1.26%"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'30.70"
$ws.Range("E3").Value = "'11.92%"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'5.170"
$ws.Range("E4").Value = "'0.37%"
$ws.Range("G4").Value = "'12"
$ws.Range("D5").Value = "'0.05729"
$ws.Range("E5").Value = "'1.56%"
$ws.Range("G5").Value = "'12"
$ws.Range("D6").Value = "'6.604"
$ws.Range("E6").Value = "'2.12%"
$ws.Range("G6").Value = "'12"
$ws.Range("D7").Value = "'3.072"
$ws.Range("E7").Value = "'2.24%"
$ws.Range("G7").Value = "'12"
$ws.Range("D8").Value = "'0.8552"
$ws.Range("E8").Value = "'4.78%"
$ws.Range("G8").Value = "'12"
$ws.Range("D9").Value = "'0.8782"
$ws.Range("E9").Value = "'5.68%"
$ws.Range("G9").Value = "'12"
$ws.Range("D10").Value = "'0.1365"
$ws.Range("E10").Value = "'2.68%"
$ws.Range("G10").Value = "'12"
$ws.Range("D11").Value = "'0.07070"
$ws.Range("E11").Value = "'2.11%"
$ws.Range("G11").Value = "'12"
$ws.Range("D12").Value = "'0.02864"
$ws.Range("E12").Value = "'-2.32%"
$ws.Range("G12").Value = "'12"
$ws.Range("D13").Value = "'0.09385"
$ws.Range("E13").Value = "'-0.08%"
$ws.Range("G13").Value = "'12"
$ws.Range("D14").Value = "'0.001513"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("G14").Value = "'12"
$ws.Range("D15").Value = "'0.04159"
$ws.Range("E15").Value = "'-2.22%"
$ws.Range("G15").Value = "'12"
$ws.Range("D16").Value = "'0.0006017"
$ws.Range("E16").Value = "'0.48%"
$ws.Range("G16").Value = "'12"
$ws.Range("D17").Value = "'0.006210"
$ws.Range("E17").Value = "'-0.10%"
$ws.Range("G17").Value = "'12"
$ws.Range("E18").Value = "'-0.47%"
$ws.Range("G18").Value = "'12"
$ws.Range("D19").Value = "'2.261"
$ws.Range("E19").Value = "'-2.16%"
$ws.Range("G19").Value = "'12"
$ws.Range("D20").Value = "'0.3162"
$ws.Range("E20").Value = "'1.58%"
$ws.Range("G20").Value = "'12"
$ws.Range("D21").Value = "'0.03233"
$ws.Range("E21").Value = "'4.33%"
$ws.Range("G21").Value = "'12"
$ws.Range("E22").Value = "'0.68%"
$ws.Range("G22").Value = "'12"
$ws.Range("D23").Value = "'2.913"
$ws.Range("E23").Value = "'-22.34%"
$ws.Range("G23").Value = "'12"
$ws.Range("E24").Value = "'0.39%"
$ws.Range("G24").Value = "'12"
$ws.Range("D25").Value = "'0.001213"
$ws.Range("E25").Value = "'-0.98%"
$ws.Range("G25").Value = "'12"
$ws.Range("D26").Value = "'0.004494"
$ws.Range("E26").Value = "'0.58%"
$ws.Range("G26").Value = "'12"
$ws.Range("E27").Value = "'23.39%"
$ws.Range("G27").Value = "'12"
$ws.Range("E28").Value = "'-0.05%"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.03783"
$ws.Range("E40").Value = "'3.73%"
$ws.Range("G40").Value = "'12"
$ws.Range("D41").Value = "'0.005685"
$ws.Range("E41").Value = "'-6.01%"
$ws.Range("G41").Value = "'12"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("E42").Value = "'1.80%"
$ws.Range("G42").Value = "'12"
$ws.Range("E43").Value = "'-4.40%"
$ws.Range("G43").Value = "'12"
$ws.Range("E44").Value = "'21.86%"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.00005084"
$ws.Range("E45").Value = "'-5.41%"
$ws.Range("G45").Value = "'12"
$ws.Range("E46").Value = "'-0.05%"
$ws.Range("G46").Value = "'12"
$ws.Range("D47").Value = "'0.07996"
$ws.Range("E47").Value = "'-20.82%"
$ws.Range("G47").Value = "'12"
$ws.Range("D48").Value = "'0.002765"
$ws.Range("E48").Value = "'-41.60%"
$ws.Range("G48").Value = "'12"
$ws.Range("E49").Value = "'-0.05%"
$ws.Range("G49").Value = "'12"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("G50").Value = "'12"
$ws.Range("G51").Value = "'12"
